$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New asset rows (row 3 through row 19) appended below the existing header
# row (row 2: name | file name | purpose | description | dimensiuons |
# poly count | texture res | status).
#
# Columns in use: B=name, C=file name, D=purpose, F=dimensiuons
# (columns E, G, H, I are left blank, same as in the header's sibling data).

$rows = @(
    @{ B = "Deer Statue";      C = "deerstatueactualfulltexture"; D = "feature asset";   F = "7x13x8" },
    @{ B = "mushroom orb";     C = "orbshroomactual";             D = "pick-ups";        F = "0.3x2" },
    @{ B = "wall";             C = "mapwallsizetest";             D = "border of map";   F = "70x30x70" },
    @{ B = "shrub";            C = "shrub1";                      D = "decoration";      F = "1x1x1 (multiple upscales)" },
    @{ B = "grass blade";      C = "grassblade";                  D = "decoration";      F = "1x1x1" },
    @{ B = "waterfall";        C = "waterfallwater";              D = "puzzle piece";    F = "11.5x11.5x11.5" },
    @{ B = "waterfall bay";    C = "waterfallbay";                D = "decoration";      F = "69x7x69" },
    @{ B = "waterfall stand";  C = "waterfallstand";              D = "puzzle piece";    F = "8.5x3.5x14" },
    @{ B = "large tree";       C = "largetree";                   D = "decoration";      F = "40x60x38" },
    @{ B = "small tree";       C = "smalltreesizetest";           D = "decoration";      F = "4x25x3.5" },
    @{ B = "rock";             C = "rock";                        D = "decoration";      F = "3x3x3" },
    @{ B = "stone";            C = "stone";                       D = "decoration";      F = "2x2x2" },
    @{ B = "pebble";           C = "pebble";                      D = "decoration";      F = "1x1x1" },
    @{ B = "ground";           D = "ground" },
    @{ B = "log";              C = "log";                         D = "decoration ";     F = "0.5x7" },
    @{ B = "stump";            C = "stump";                       D = "decoration";      F = "0.5x2" },
    @{ B = "pathway";          D = "path" }
)

$r = 3
foreach ($row in $rows) {
    if ($row.ContainsKey("B")) { $ws.Range("B$r").Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Range("C$r").Value = $row.C }
    if ($row.ContainsKey("D")) { $ws.Range("D$r").Value = $row.D }
    if ($row.ContainsKey("F")) { $ws.Range("F$r").Value = $row.F }
    $r++
}

# Column widths widened to fit the newly-entered data.
$ws.Columns.Item(2).ColumnWidth = 13.73
$ws.Columns.Item(3).ColumnWidth = 26.63
$ws.Columns.Item(4).ColumnWidth = 14.18
$ws.Columns.Item(5).ColumnWidth = 18.45
$ws.Columns.Item(6).ColumnWidth = 22.36
$ws.Columns.Item(7).ColumnWidth = 11.18
$ws.Columns.Item(8).ColumnWidth = 11.45
$ws.Columns.Item(9).ColumnWidth = 9.18

# Selection ends on F8, matching where editing left off.
$ws.Range("F8").Select()
